$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(5, 2).Value = 6430240
$ws.Cells.Item(5, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(5, 6).Value = 'LASK Linz'
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 2
$ws.Cells.Item(5, 12).Value = 1.8
$ws.Cells.Item(5, 13).Value = 4
$ws.Cells.Item(5, 14).Value = 3.8
$ws.Cells.Item(5, 15).Value = 1.727
$ws.Cells.Item(5, 16).Value = 4.333
$ws.Cells.Item(5, 17).Value = 4
$ws.Cells.Item(5, 18).Value = -0.75
$ws.Cells.Item(5, 19).Value = 1.925
$ws.Cells.Item(5, 20).Value = 1.925
$ws.Cells.Item(5, 21).Value = 3.25
$ws.Cells.Item(5, 22).Value = 1.975
$ws.Cells.Item(5, 23).Value = 1.875
$ws.Cells.Item(5, 24).Value = 0.7270000000000001
$ws.Cells.Item(5, 27).Value = 0.925
$ws.Cells.Item(5, 29).Value = -1
$ws.Cells.Item(5, 30).Value = 0.875
$ws.Cells.Item(7, 2).Value = 6430241
$ws.Cells.Item(7, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(7, 6).Value = 'Rapid Vienna'
$ws.Cells.Item(7, 8).Value = 1
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 12).Value = 4.2
$ws.Cells.Item(7, 13).Value = 3.75
$ws.Cells.Item(7, 14).Value = 1.8
$ws.Cells.Item(7, 15).Value = 4.2
$ws.Cells.Item(7, 16).Value = 4
$ws.Cells.Item(7, 17).Value = 1.8
$ws.Cells.Item(7, 18).Value = 0.75
$ws.Cells.Item(7, 19).Value = 1.85
$ws.Cells.Item(7, 20).Value = 2
$ws.Cells.Item(7, 21).Value = 3
$ws.Cells.Item(7, 22).Value = 2.025
$ws.Cells.Item(7, 23).Value = 1.825
$ws.Cells.Item(7, 24).Value = 3.2
$ws.Cells.Item(7, 27).Value = 0.8500000000000001
$ws.Cells.Item(7, 29).Value = 0
$ws.Cells.Item(7, 30).Value = 0
$ws.Cells.Item(10, 2).Value = 6851964
$ws.Cells.Item(10, 5).Value = 'Wolfsberger AC'
$ws.Cells.Item(10, 6).Value = 'FC Blau Weiss Linz'
$ws.Cells.Item(10, 7).Value = 2
$ws.Cells.Item(10, 8).Value = 1
$ws.Cells.Item(10, 9).Value = 1
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 'H'
$ws.Cells.Item(10, 12).Value = 1.909
$ws.Cells.Item(10, 13).Value = 3.5
$ws.Cells.Item(10, 14).Value = 3.5
$ws.Cells.Item(10, 15).Value = 2.2
$ws.Cells.Item(10, 16).Value = 3.4
$ws.Cells.Item(10, 17).Value = 3.25
$ws.Cells.Item(10, 19).Value = 1.975
$ws.Cells.Item(10, 20).Value = 1.875
$ws.Cells.Item(10, 21).Value = 2.5
$ws.Cells.Item(10, 22).Value = 2
$ws.Cells.Item(10, 23).Value = 1.85
$ws.Cells.Item(10, 24).Value = 1.2
$ws.Cells.Item(10, 26).Value = -1
$ws.Cells.Item(10, 27).Value = 0.9750000000000001
$ws.Cells.Item(10, 28).Value = -1
$ws.Cells.Item(10, 29).Value = 1
$ws.Cells.Item(11, 2).Value = 6847027
$ws.Cells.Item(11, 5).Value = 'WSG Swarovski Tirol'
$ws.Cells.Item(11, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(11, 7).Value = 1
$ws.Cells.Item(11, 8).Value = 3
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 1
$ws.Cells.Item(11, 11).Value = 'A'
$ws.Cells.Item(11, 12).Value = 2.1
$ws.Cells.Item(11, 13).Value = 3.4
$ws.Cells.Item(11, 14).Value = 3.1
$ws.Cells.Item(11, 15).Value = 2.3
$ws.Cells.Item(11, 16).Value = 3.6
$ws.Cells.Item(11, 17).Value = 2.8
$ws.Cells.Item(11, 19).Value = 2.025
$ws.Cells.Item(11, 20).Value = 1.825
$ws.Cells.Item(11, 21).Value = 2.75
$ws.Cells.Item(11, 22).Value = 1.825
$ws.Cells.Item(11, 23).Value = 2.025
$ws.Cells.Item(11, 24).Value = -1
$ws.Cells.Item(11, 26).Value = 1.8
$ws.Cells.Item(11, 27).Value = -1
$ws.Cells.Item(11, 28).Value = 0.825
$ws.Cells.Item(11, 29).Value = 0.825
$ws.Cells.Item(13, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(15, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(19, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(20, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(20, 6).Value = 'SK Sturm Graz'
$ws.Cells.Item(27, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(28, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(32, 2).Value = 6851960
$ws.Cells.Item(32, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(32, 6).Value = 'FC Blau Weiss Linz'
$ws.Cells.Item(32, 7).Value = 4
$ws.Cells.Item(32, 8).Value = 1
$ws.Cells.Item(32, 9).Value = 2
$ws.Cells.Item(32, 11).Value = 'H'
$ws.Cells.Item(32, 12).Value = 1.285
$ws.Cells.Item(32, 13).Value = 5
$ws.Cells.Item(32, 14).Value = 7.5
$ws.Cells.Item(32, 15).Value = 1.333
$ws.Cells.Item(32, 16).Value = 5.75
$ws.Cells.Item(32, 17).Value = 7.5
$ws.Cells.Item(32, 18).Value = -1.5
$ws.Cells.Item(32, 21).Value = 3
$ws.Cells.Item(32, 22).Value = 1.975
$ws.Cells.Item(32, 23).Value = 1.875
$ws.Cells.Item(32, 24).Value = 0.333
$ws.Cells.Item(32, 26).Value = -1
$ws.Cells.Item(32, 27).Value = 0.95
$ws.Cells.Item(32, 28).Value = -1
$ws.Cells.Item(32, 29).Value = 0.9750000000000001
$ws.Cells.Item(32, 30).Value = -1
$ws.Cells.Item(33, 2).Value = 6846462
$ws.Cells.Item(33, 5).Value = 'Wolfsberger AC'
$ws.Cells.Item(33, 6).Value = 'FC Salzburg'
$ws.Cells.Item(33, 7).Value = 1
$ws.Cells.Item(33, 8).Value = 2
$ws.Cells.Item(33, 9).Value = 1
$ws.Cells.Item(33, 11).Value = 'A'
$ws.Cells.Item(33, 12).Value = 6.5
$ws.Cells.Item(33, 13).Value = 5.1
$ws.Cells.Item(33, 14).Value = 1.3
$ws.Cells.Item(33, 15).Value = 5.75
$ws.Cells.Item(33, 16).Value = 4.2
$ws.Cells.Item(33, 17).Value = 1.533
$ws.Cells.Item(33, 18).Value = 1
$ws.Cells.Item(33, 21).Value = 2.75
$ws.Cells.Item(33, 22).Value = 1.875
$ws.Cells.Item(33, 23).Value = 1.975
$ws.Cells.Item(33, 24).Value = -1
$ws.Cells.Item(33, 26).Value = 0.5329999999999999
$ws.Cells.Item(33, 27).Value = 0
$ws.Cells.Item(33, 28).Value = 0
$ws.Cells.Item(33, 29).Value = 0.4375
$ws.Cells.Item(33, 30).Value = -0.5
$ws.Cells.Item(36, 6).Value = 'SK Sturm Graz'
$ws.Cells.Item(38, 2).Value = 6847045
$ws.Cells.Item(38, 5).Value = 'Wolfsberger AC'
$ws.Cells.Item(38, 6).Value = 'Hartberg'
$ws.Cells.Item(38, 7).Value = 0
$ws.Cells.Item(38, 8).Value = 3
$ws.Cells.Item(38, 9).Value = 0
$ws.Cells.Item(38, 10).Value = 1
$ws.Cells.Item(38, 12).Value = 2.05
$ws.Cells.Item(38, 13).Value = 3.3
$ws.Cells.Item(38, 14).Value = 3.1
$ws.Cells.Item(38, 15).Value = 2.3
$ws.Cells.Item(38, 16).Value = 3.3
$ws.Cells.Item(38, 17).Value = 3
$ws.Cells.Item(38, 18).Value = -0.25
$ws.Cells.Item(38, 21).Value = 2.5
$ws.Cells.Item(38, 26).Value = 2
$ws.Cells.Item(39, 2).Value = 6851959
$ws.Cells.Item(39, 5).Value = 'WSG Swarovski Tirol'
$ws.Cells.Item(39, 6).Value = 'FC Blau Weiss Linz'
$ws.Cells.Item(39, 7).Value = 2
$ws.Cells.Item(39, 8).Value = 4
$ws.Cells.Item(39, 9).Value = 1
$ws.Cells.Item(39, 10).Value = 2
$ws.Cells.Item(39, 12).Value = 2.1
$ws.Cells.Item(39, 13).Value = 3.4
$ws.Cells.Item(39, 14).Value = 2.9
$ws.Cells.Item(39, 15).Value = 2.7
$ws.Cells.Item(39, 16).Value = 3.75
$ws.Cells.Item(39, 17).Value = 2.5
$ws.Cells.Item(39, 18).Value = 0
$ws.Cells.Item(39, 21).Value = 2.75
$ws.Cells.Item(39, 26).Value = 1.5
$ws.Cells.Item(40, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(43, 6).Value = 'SK Sturm Graz'
$ws.Cells.Item(44, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(46, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(51, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(55, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(57, 6).Value = 'SK Sturm Graz'
$ws.Cells.Item(59, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(65, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(67, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(68, 2).Value = 6846467
$ws.Cells.Item(68, 5).Value = 'FC Salzburg'
$ws.Cells.Item(68, 6).Value = 'LASK Linz'
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 11).Value = 'A'
$ws.Cells.Item(68, 12).Value = 1.533
$ws.Cells.Item(68, 13).Value = 3.8
$ws.Cells.Item(68, 15).Value = 1.55
$ws.Cells.Item(68, 17).Value = 6
$ws.Cells.Item(68, 18).Value = -1
$ws.Cells.Item(68, 21).Value = 3
$ws.Cells.Item(68, 22).Value = 2.05
$ws.Cells.Item(68, 23).Value = 1.8
$ws.Cells.Item(68, 24).Value = -1
$ws.Cells.Item(68, 26).Value = 5
$ws.Cells.Item(68, 27).Value = -1
$ws.Cells.Item(68, 28).Value = 0.925
$ws.Cells.Item(68, 29).Value = -1
$ws.Cells.Item(68, 30).Value = 0.8
$ws.Cells.Item(69, 2).Value = 6847066
$ws.Cells.Item(69, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(69, 6).Value = 'Hartberg'
$ws.Cells.Item(69, 7).Value = 2
$ws.Cells.Item(69, 9).Value = 2
$ws.Cells.Item(69, 11).Value = 'H'
$ws.Cells.Item(69, 12).Value = 1.5
$ws.Cells.Item(69, 13).Value = 4
$ws.Cells.Item(69, 15).Value = 1.7
$ws.Cells.Item(69, 17).Value = 4.75
$ws.Cells.Item(69, 18).Value = -0.75
$ws.Cells.Item(69, 21).Value = 2.5
$ws.Cells.Item(69, 22).Value = 1.825
$ws.Cells.Item(69, 23).Value = 2.025
$ws.Cells.Item(69, 24).Value = 0.7
$ws.Cells.Item(69, 26).Value = -1
$ws.Cells.Item(69, 27).Value = 0.4625
$ws.Cells.Item(69, 28).Value = -0.5
$ws.Cells.Item(69, 29).Value = 0.825
$ws.Cells.Item(69, 30).Value = -1
$ws.Cells.Item(73, 6).Value = 'SK Sturm Graz'
$ws.Cells.Item(75, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(77, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(80, 6).Value = 'SK Sturm Graz'
$ws.Cells.Item(85, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(89, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(89, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(94, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(96, 6).Value = 'SK Sturm Graz'
$ws.Cells.Item(98, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(101, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(106, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(107, 2).Value = 7529184
$ws.Cells.Item(107, 5).Value = 'FC Blau Weiss Linz'
$ws.Cells.Item(107, 6).Value = 'WSG Swarovski Tirol'
$ws.Cells.Item(107, 7).Value = 1
$ws.Cells.Item(107, 8).Value = 2
$ws.Cells.Item(107, 10).Value = 1
$ws.Cells.Item(107, 11).Value = 'A'
$ws.Cells.Item(107, 12).Value = 1.909
$ws.Cells.Item(107, 13).Value = 3.4
$ws.Cells.Item(107, 14).Value = 3.5
$ws.Cells.Item(107, 15).Value = 2.15
$ws.Cells.Item(107, 16).Value = 3.4
$ws.Cells.Item(107, 17).Value = 3.3
$ws.Cells.Item(107, 18).Value = -0.25
$ws.Cells.Item(107, 19).Value = 1.85
$ws.Cells.Item(107, 20).Value = 2
$ws.Cells.Item(107, 21).Value = 2.5
$ws.Cells.Item(107, 22).Value = 1.95
$ws.Cells.Item(107, 23).Value = 1.9
$ws.Cells.Item(107, 24).Value = -1
$ws.Cells.Item(107, 26).Value = 2.3
$ws.Cells.Item(107, 27).Value = -1
$ws.Cells.Item(107, 28).Value = 1
$ws.Cells.Item(107, 29).Value = 0.95
$ws.Cells.Item(107, 30).Value = -1
$ws.Cells.Item(108, 2).Value = 6847093
$ws.Cells.Item(108, 5).Value = 'Hartberg'
$ws.Cells.Item(108, 6).Value = 'Wolfsberger AC'
$ws.Cells.Item(108, 7).Value = 2
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 11).Value = 'H'
$ws.Cells.Item(108, 12).Value = 2.05
$ws.Cells.Item(108, 13).Value = 3.5
$ws.Cells.Item(108, 14).Value = 3
$ws.Cells.Item(108, 15).Value = 2.55
$ws.Cells.Item(108, 16).Value = 3.5
$ws.Cells.Item(108, 17).Value = 2.625
$ws.Cells.Item(108, 18).Value = 0
$ws.Cells.Item(108, 19).Value = 1.925
$ws.Cells.Item(108, 20).Value = 1.925
$ws.Cells.Item(108, 21).Value = 2.75
$ws.Cells.Item(108, 22).Value = 2.025
$ws.Cells.Item(108, 23).Value = 1.825
$ws.Cells.Item(108, 24).Value = 1.55
$ws.Cells.Item(108, 26).Value = -1
$ws.Cells.Item(108, 27).Value = 0.925
$ws.Cells.Item(108, 28).Value = -1
$ws.Cells.Item(108, 29).Value = -1
$ws.Cells.Item(108, 30).Value = 0.825
$ws.Cells.Item(109, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(110, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(115, 6).Value = 'SK Sturm Graz'
$ws.Cells.Item(120, 6).Value = 'SK Sturm Graz'
$ws.Cells.Item(121, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(125, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(126, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(128, 6).Value = 'SK Sturm Graz'
$ws.Cells.Item(130, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(136, 2).Value = 6847111
$ws.Cells.Item(136, 5).Value = 'Wolfsberger AC'
$ws.Cells.Item(136, 6).Value = 'SCR Altach'
$ws.Cells.Item(136, 9).Value = 1
$ws.Cells.Item(136, 12).Value = 2.05
$ws.Cells.Item(136, 13).Value = 3.5
$ws.Cells.Item(136, 14).Value = 3.5
$ws.Cells.Item(136, 15).Value = 2.625
$ws.Cells.Item(136, 16).Value = 3.2
$ws.Cells.Item(136, 17).Value = 2.75
$ws.Cells.Item(136, 18).Value = 0
$ws.Cells.Item(136, 19).Value = 1.875
$ws.Cells.Item(136, 20).Value = 1.975
$ws.Cells.Item(136, 21).Value = 2
$ws.Cells.Item(136, 22).Value = 1.75
$ws.Cells.Item(136, 23).Value = 2.05
$ws.Cells.Item(136, 25).Value = 2.2
$ws.Cells.Item(136, 27).Value = 0
$ws.Cells.Item(136, 28).Value = 0
$ws.Cells.Item(136, 29).Value = 0
$ws.Cells.Item(136, 30).Value = 0
$ws.Cells.Item(137, 2).Value = 6847112
$ws.Cells.Item(137, 5).Value = 'FK Austria Vienna'
$ws.Cells.Item(137, 6).Value = 'WSG Swarovski Tirol'
$ws.Cells.Item(137, 7).Value = 2
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 9).Value = 1
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 'H'
$ws.Cells.Item(137, 12).Value = 1.5
$ws.Cells.Item(137, 13).Value = 4.333
$ws.Cells.Item(137, 14).Value = 6
$ws.Cells.Item(137, 15).Value = 1.5
$ws.Cells.Item(137, 16).Value = 4.5
$ws.Cells.Item(137, 17).Value = 6.5
$ws.Cells.Item(137, 18).Value = -1.25
$ws.Cells.Item(137, 19).Value = 2.05
$ws.Cells.Item(137, 20).Value = 1.8
$ws.Cells.Item(137, 21).Value = 2.75
$ws.Cells.Item(137, 22).Value = 1.9
$ws.Cells.Item(137, 23).Value = 1.95
$ws.Cells.Item(137, 24).Value = 0.5
$ws.Cells.Item(137, 25).Value = -1
$ws.Cells.Item(137, 27).Value = 1.05
$ws.Cells.Item(137, 28).Value = -1
$ws.Cells.Item(137, 29).Value = -1
$ws.Cells.Item(137, 30).Value = 0.95
$ws.Cells.Item(138, 2).Value = 6847113
$ws.Cells.Item(138, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(138, 6).Value = 'Rapid Vienna'
$ws.Cells.Item(138, 7).Value = 1
$ws.Cells.Item(138, 8).Value = 1
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 1
$ws.Cells.Item(138, 11).Value = 'D'
$ws.Cells.Item(138, 12).Value = 3.6
$ws.Cells.Item(138, 13).Value = 3.8
$ws.Cells.Item(138, 14).Value = 1.909
$ws.Cells.Item(138, 15).Value = 3.6
$ws.Cells.Item(138, 16).Value = 2.75
$ws.Cells.Item(138, 17).Value = 2.375
$ws.Cells.Item(138, 18).Value = 0.25
$ws.Cells.Item(138, 19).Value = 1.85
$ws.Cells.Item(138, 20).Value = 2
$ws.Cells.Item(138, 21).Value = 2.25
$ws.Cells.Item(138, 22).Value = 1.875
$ws.Cells.Item(138, 23).Value = 1.975
$ws.Cells.Item(138, 24).Value = -1
$ws.Cells.Item(138, 25).Value = 1.75
$ws.Cells.Item(138, 27).Value = 0.425
$ws.Cells.Item(138, 28).Value = -0.5
$ws.Cells.Item(138, 29).Value = -0.5
$ws.Cells.Item(138, 30).Value = 0.4875
$ws.Cells.Item(139, 2).Value = 6847114
$ws.Cells.Item(139, 5).Value = 'Hartberg'
$ws.Cells.Item(139, 6).Value = 'SK Sturm Graz'
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 12).Value = 3.4
$ws.Cells.Item(139, 13).Value = 3.6
$ws.Cells.Item(139, 14).Value = 2
$ws.Cells.Item(139, 15).Value = 3.3
$ws.Cells.Item(139, 16).Value = 3.3
$ws.Cells.Item(139, 17).Value = 2.2
$ws.Cells.Item(139, 18).Value = 0.25
$ws.Cells.Item(139, 19).Value = 1.925
$ws.Cells.Item(139, 20).Value = 1.925
$ws.Cells.Item(139, 21).Value = 2.25
$ws.Cells.Item(139, 22).Value = 1.9
$ws.Cells.Item(139, 23).Value = 1.95
$ws.Cells.Item(139, 25).Value = 2.3
$ws.Cells.Item(139, 27).Value = 0.4625
$ws.Cells.Item(139, 28).Value = -0.5
$ws.Cells.Item(139, 29).Value = -0.5
$ws.Cells.Item(139, 30).Value = 0.475
$ws.Cells.Item(144, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(144, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(150, 6).Value = 'SK Sturm Graz'
$ws.Cells.Item(151, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(155, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(156, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(162, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(163, 6).Value = 'SK Sturm Graz'
$ws.Cells.Item(164, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(169, 6).Value = 'SK Sturm Graz'
$ws.Cells.Item(173, 2).Value = 7948256
$ws.Cells.Item(173, 5).Value = 'Hartberg'
$ws.Cells.Item(173, 6).Value = 'LASK Linz'
$ws.Cells.Item(173, 7).Value = 1
$ws.Cells.Item(173, 8).Value = 2
$ws.Cells.Item(173, 9).Value = 1
$ws.Cells.Item(173, 10).Value = 1
$ws.Cells.Item(173, 11).Value = 'A'
$ws.Cells.Item(173, 12).Value = 3
$ws.Cells.Item(173, 13).Value = 3.5
$ws.Cells.Item(173, 14).Value = 2.25
$ws.Cells.Item(173, 15).Value = 3.2
$ws.Cells.Item(173, 16).Value = 3.5
$ws.Cells.Item(173, 17).Value = 2.2
$ws.Cells.Item(173, 18).Value = 0.25
$ws.Cells.Item(173, 19).Value = 1.975
$ws.Cells.Item(173, 20).Value = 1.875
$ws.Cells.Item(173, 21).Value = 2.25
$ws.Cells.Item(173, 22).Value = 1.8
$ws.Cells.Item(173, 23).Value = 2.05
$ws.Cells.Item(173, 24).Value = -1
$ws.Cells.Item(173, 26).Value = 1.2
$ws.Cells.Item(173, 27).Value = -1
$ws.Cells.Item(173, 28).Value = 0.875
$ws.Cells.Item(173, 29).Value = 0.8
$ws.Cells.Item(174, 2).Value = 7948255
$ws.Cells.Item(174, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(174, 6).Value = 'FC Salzburg'
$ws.Cells.Item(174, 7).Value = 4
$ws.Cells.Item(174, 8).Value = 3
$ws.Cells.Item(174, 9).Value = 0
$ws.Cells.Item(174, 10).Value = 2
$ws.Cells.Item(174, 11).Value = 'H'
$ws.Cells.Item(174, 12).Value = 6.5
$ws.Cells.Item(174, 13).Value = 4.5
$ws.Cells.Item(174, 14).Value = 1.444
$ws.Cells.Item(174, 15).Value = 9
$ws.Cells.Item(174, 16).Value = 5.25
$ws.Cells.Item(174, 17).Value = 1.333
$ws.Cells.Item(174, 18).Value = 1.5
$ws.Cells.Item(174, 19).Value = 1.925
$ws.Cells.Item(174, 20).Value = 1.925
$ws.Cells.Item(174, 21).Value = 3
$ws.Cells.Item(174, 22).Value = 1.825
$ws.Cells.Item(174, 23).Value = 2.025
$ws.Cells.Item(174, 24).Value = 8
$ws.Cells.Item(174, 26).Value = -1
$ws.Cells.Item(174, 27).Value = 0.925
$ws.Cells.Item(174, 28).Value = -1
$ws.Cells.Item(174, 29).Value = 0.825
$ws.Cells.Item(175, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(179, 2).Value = 7948258
$ws.Cells.Item(179, 5).Value = 'Hartberg'
$ws.Cells.Item(179, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(179, 7).Value = 3
$ws.Cells.Item(179, 8).Value = 2
$ws.Cells.Item(179, 12).Value = 2.3
$ws.Cells.Item(179, 13).Value = 3.3
$ws.Cells.Item(179, 14).Value = 3.1
$ws.Cells.Item(179, 17).Value = 3.3
$ws.Cells.Item(179, 19).Value = 1.8
$ws.Cells.Item(179, 20).Value = 2.05
$ws.Cells.Item(179, 21).Value = 2.75
$ws.Cells.Item(179, 27).Value = 0.8
$ws.Cells.Item(180, 2).Value = 7948257
$ws.Cells.Item(180, 5).Value = 'LASK Linz'
$ws.Cells.Item(180, 6).Value = 'Rapid Vienna'
$ws.Cells.Item(180, 7).Value = 5
$ws.Cells.Item(180, 8).Value = 0
$ws.Cells.Item(180, 12).Value = 1.5
$ws.Cells.Item(180, 13).Value = 3
$ws.Cells.Item(180, 14).Value = 1.5
$ws.Cells.Item(180, 17).Value = 3.5
$ws.Cells.Item(180, 19).Value = 1.825
$ws.Cells.Item(180, 20).Value = 2.025
$ws.Cells.Item(180, 21).Value = 2.5
$ws.Cells.Item(180, 27).Value = 0.825
$ws.Cells.Item(181, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(185, 5).Value = 'Austria Klagenfurt'
$ws.Cells.Item(186, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(192, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(193, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(197, 2).Value = 7947241
$ws.Cells.Item(197, 5).Value = 'SK Sturm Graz'
$ws.Cells.Item(197, 6).Value = 'Austria Klagenfurt'
$ws.Cells.Item(197, 7).Value = 2
$ws.Cells.Item(197, 8).Value = 0
$ws.Cells.Item(197, 9).Value = 0
$ws.Cells.Item(197, 10).Value = 0
$ws.Cells.Item(197, 12).Value = 1.4
$ws.Cells.Item(197, 13).Value = 5
$ws.Cells.Item(197, 14).Value = 7
$ws.Cells.Item(197, 15).Value = 1.333
$ws.Cells.Item(197, 16).Value = 5.75
$ws.Cells.Item(197, 17).Value = 8
$ws.Cells.Item(197, 19).Value = 1.85
$ws.Cells.Item(197, 20).Value = 2
$ws.Cells.Item(197, 21).Value = 3.25
$ws.Cells.Item(197, 22).Value = 1.975
$ws.Cells.Item(197, 23).Value = 1.875
$ws.Cells.Item(197, 24).Value = 0.333
$ws.Cells.Item(197, 27).Value = 0.8500000000000001
$ws.Cells.Item(197, 29).Value = -1
$ws.Cells.Item(197, 30).Value = 0.875
$ws.Cells.Item(198, 2).Value = 7948263
$ws.Cells.Item(198, 5).Value = 'FC Salzburg'
$ws.Cells.Item(198, 6).Value = 'LASK Linz'
$ws.Cells.Item(198, 7).Value = 7
$ws.Cells.Item(198, 8).Value = 1
$ws.Cells.Item(198, 9).Value = 3
$ws.Cells.Item(198, 10).Value = 1
$ws.Cells.Item(198, 12).Value = 1.55
$ws.Cells.Item(198, 13).Value = 4.333
$ws.Cells.Item(198, 14).Value = 5.25
$ws.Cells.Item(198, 15).Value = 1.4
$ws.Cells.Item(198, 16).Value = 5.25
$ws.Cells.Item(198, 17).Value = 6.5
$ws.Cells.Item(198, 19).Value = 2.025
$ws.Cells.Item(198, 20).Value = 1.825
$ws.Cells.Item(198, 21).Value = 3.5
$ws.Cells.Item(198, 22).Value = 1.95
$ws.Cells.Item(198, 23).Value = 1.9
$ws.Cells.Item(198, 24).Value = 0.3999999999999999
$ws.Cells.Item(198, 27).Value = 1.025
$ws.Cells.Item(198, 29).Value = 0.95
$ws.Cells.Item(198, 30).Value = -1
